# Journal de travaille - add new entry row (Mise au propre, 2021-03-12)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the table ("Tableau1") by one row; this keeps the ListObject's
# range / AutoFilter / sheet dimension all in sync with the new row.
$lo = $ws.ListObjects.Item("Tableau1")
$lo.ListRows.Add() | Out-Null

# --- New row (row 44) content -------------------------------------------
$ws.Range("B44").Value = 44267
$ws.Range("C44").Value = 0.57291666666666663
$ws.Range("D44").Value = 0.59375
$ws.Range("E44").Formula = "=IF(ISBLANK(Tableau1[[#This Row],[Heure fin]]),`"`",Tableau1[[#This Row],[Heure fin]]-Tableau1[[#This Row],[Heure début]])"
$ws.Range("F44").Value = "Ma-20"
$ws.Range("G44").Value = "Code"
$ws.Range("H44").Value = "Mise au propre"
$ws.Range("I44").Value = "CPNV"
$ws.Range("J44").Value = "J'ai corrigé d es fautes d'orthgraphe, ajouter des legends et aussi bien placer mes bateau"
$ws.Range("K44").Value = "Oui"

# Match the date's number format (mm-dd-yy / style "1") by copying the
# format from the cell above it, same as Excel does when a table row is
# filled in by hand.
$ws.Range("B43").Copy() | Out-Null
$ws.Range("B44").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# The row wraps to three lines at this column width, same as the other
# multi-line "Descriptif" rows in the sheet.
$ws.Rows.Item(44).RowHeight = 43.2

# Selection ends on K44 after the edit.
$ws.Range("K44").Select() | Out-Null
